# Endless mode - post-level popup draft version
# Adds a "Dynamic music" story row, a new "Endless runner mode" section with
# several stories on the Stories sheet, and three supporting rows on the
# Content sheet (Jungle of Tikirara section).

$wb = $excel.ActiveWorkbook
$wsStories = $wb.Worksheets.Item("Stories")
$wsContent = $wb.Worksheets.Item("Content")

# ---------------------------------------------------------------------
# Styles: column F becomes a left-aligned, wrapped column (wider too).
# ---------------------------------------------------------------------
$wsStories.Columns.Item(6).ColumnWidth = 31.17
$wsStories.Range("F1").HorizontalAlignment = -4131
$wsStories.Range("F1").WrapText = $true

# ---------------------------------------------------------------------
# Stories sheet: insert the "Dynamic music" row above row 12.
# ---------------------------------------------------------------------
$wsStories.Rows.Item(12).Insert()

$wsStories.Range("A12").Value = "Dynamic music"
$wsStories.Range("B12").Value = "Low"
$wsStories.Range("C12").Value = "1.5"
$wsStories.Range("F12").Value = "Play combat music when there are enemies around"
$wsStories.Range("F12").HorizontalAlignment = -4131
$wsStories.Range("F12").WrapText = $true
$wsStories.Rows.Item(12).RowHeight = 30

# Existing F-column notes now pick up the left-aligned wrap style too.
$wsStories.Range("F13").HorizontalAlignment = -4131
$wsStories.Range("F13").WrapText = $true
$wsStories.Range("F14").HorizontalAlignment = -4131
$wsStories.Range("F14").WrapText = $true
$wsStories.Range("F41").HorizontalAlignment = -4131
$wsStories.Range("F41").WrapText = $true
$wsStories.Range("F42").HorizontalAlignment = -4131
$wsStories.Range("F42").WrapText = $true
$wsStories.Range("F50").HorizontalAlignment = -4131
$wsStories.Range("F50").WrapText = $true

# ---------------------------------------------------------------------
# Content sheet: "Macaw" note under the Jungle of Tikirara section.
# ---------------------------------------------------------------------
$wsContent.Range("A23").Value = "Macaw"

# ---------------------------------------------------------------------
# Stories sheet: new "Endless runner mode" section (rows 53-55).
# ---------------------------------------------------------------------
$wsStories.Range("A53").Value = "Endless runner mode"
$wsStories.Range("A53").Font.Bold = $true

$wsStories.Range("A54").Value = "Make tiles inactive after t-junction"

$wsStories.Range("A55").Value = "Keep only one tile per type in recycle list and try destroying others"
$wsStories.Rows.Item(55).RowHeight = 30

# ---------------------------------------------------------------------
# Content sheet: two more stumble/obstacle notes.
# ---------------------------------------------------------------------
$wsContent.Range("A24").Value = "Add stumbles: roots and rocks"
$wsContent.Range("A25").Value = "Add objects in middlle of the path"

# ---------------------------------------------------------------------
# Stories sheet: remaining "Endless runner mode" stories (rows 56-66).
# ---------------------------------------------------------------------
$wsStories.Range("A56").Value = "Add a way to tag tiles as suitable for endless running with a weight system"
$wsStories.Range("F56").Value = "isLandmark, isGoodForEndless, probability of apperance"
$wsStories.Range("F56").HorizontalAlignment = -4131
$wsStories.Range("F56").WrapText = $true
$wsStories.Rows.Item(56).RowHeight = 30

$wsStories.Range("A57").Value = "New personal high score message in-game"

$wsStories.Range("A58").Value = "Have a high score per theme"

$wsStories.Range("A59").Value = "Have a leaderboard on gamecenter per theme"
$wsStories.Rows.Item(59).RowHeight = 30

$wsStories.Range("A60").Value = "Have a challenge a friend mode"
$wsStories.Range("F60").Value = "also challenge non-app users"
$wsStories.Range("F60").HorizontalAlignment = -4131
$wsStories.Range("F60").WrapText = $true

$wsStories.Range("A61").Value = "Code Message center - display challenges and lives received"
$wsStories.Rows.Item(61).RowHeight = 30

$wsStories.Range("A62").Value = "Code new pre-level popup"

$wsStories.Range("A63").Value = "Code new post-level popup"

$wsStories.Range("A64").Value = "Code new save me popup"

$wsStories.Range("A65").Value = "Code selecting theme for endless running"
$wsStories.Range("F65").Value = "Maybe as a scrollable bar"
$wsStories.Range("F65").HorizontalAlignment = -4131
$wsStories.Range("F65").WrapText = $true

$wsStories.Range("A66").Value = "Inform user when he has completed an episode, that it is now unlocked for endless running"
$wsStories.Rows.Item(66).RowHeight = 45

# ---------------------------------------------------------------------
# Window / selection state.
# ---------------------------------------------------------------------
$wsContent.Activate()
$wsContent.Range("C27").Select()
$excel.ActiveWindow.ScrollRow = 2

$wsStories.Activate()
$wsStories.Range("A67").Select()
$excel.ActiveWindow.ScrollRow = 47
